$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 23611110
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H74").Value = 9900.4
$ws.Range("I74").Value = 9900.4
$ws.Range("K74").Value = 9900.4
$ws.Range("M74").Value = -8964.4
$ws.Range("H77").Value = 9900.4
$ws.Range("I77").Value = 9900.4
$ws.Range("K77").Value = 49502
$ws.Range("M77").Value = -44822
$ws.Range("H103").Value = 1675.2858
$ws.Range("J103").Value = 1621.5
$ws.Range("L103").Value = 4864.5
$ws.Range("N103").Value = -6036.5
$ws.Range("H137").Value = 1793485.5
$ws.Range("I137").Value = 1133.762
$ws.Range("J137").Value = 5557424
$ws.Range("K137").Value = 3401.286
$ws.Range("L137").Value = 16672272
$ws.Range("M137").Value = -851.2860000000001
$ws.Range("N137").Value = -16677372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1407.7097
$ws.Range("I2").Value = 1493.579
$ws.Range("J2").Value = 1271.75
$ws.Range("K2").Value = 1493.579
$ws.Range("L2").Value = 1271.75
$ws.Range("M2").Value = -1380.579
$ws.Range("N2").Value = -1497.75
$ws.Range("H4").Value = 202.2
$ws.Range("I4").Value = 217.75
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 217.75
$ws.Range("L4").Value = 140
$ws.Range("M4").Value = -101.75
$ws.Range("N4").Value = -372
$ws.Range("H5").Value = 293.4
$ws.Range("I5").Value = 293.4
$ws.Range("K5").Value = 293.4
$ws.Range("M5").Value = -181.4
$ws.Range("H32").Value = 23025908
$ws.Range("I32").Value = 27122670
$ws.Range("K32").Value = 27122670
$ws.Range("M32").Value = -27122383
$ws.Range("H74").Value = 2389.3022
$ws.Range("I74").Value = 2098.361
$ws.Range("J74").Value = 3885.5715
$ws.Range("K74").Value = 2098.361
$ws.Range("L74").Value = 3885.5715
$ws.Range("M74").Value = -1224.361
$ws.Range("N74").Value = -5633.5715
$ws.Range("H77").Value = 2389.3022
$ws.Range("I77").Value = 2098.361
$ws.Range("J77").Value = 3885.5715
$ws.Range("K77").Value = 10491.805
$ws.Range("L77").Value = 19427.8575
$ws.Range("M77").Value = -6123.805
$ws.Range("N77").Value = -28163.8575
$ws.Range("H116").Value = 1407.7097
$ws.Range("I116").Value = 1493.579
$ws.Range("J116").Value = 1271.75
$ws.Range("K116").Value = 1493.579
$ws.Range("L116").Value = 1271.75
$ws.Range("M116").Value = 800.421
$ws.Range("N116").Value = -5859.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1407.7097
$ws.Range("I3").Value = 1493.579
$ws.Range("J3").Value = 1271.75
$ws.Range("K3").Value = 1493.579
$ws.Range("L3").Value = 1271.75
$ws.Range("M3").Value = -1379.579
$ws.Range("N3").Value = -1499.75
$ws.Range("H4").Value = 293.4
$ws.Range("I4").Value = 293.4
$ws.Range("K4").Value = 293.4
$ws.Range("M4").Value = -178.4
$ws.Range("H94").Value = 946.0526
$ws.Range("I94").Value = 1056.75
$ws.Range("K94").Value = 1056.75
$ws.Range("M94").Value = -605.75
$ws.Range("H107").Value = 1374.75
$ws.Range("J107").Value = 1426.7142
$ws.Range("L107").Value = 1426.7142
$ws.Range("N107").Value = -5266.7142
$ws.Range("H134").Value = 5955388.5
$ws.Range("I134").Value = 7145016
$ws.Range("K134").Value = 21435048
$ws.Range("M134").Value = -21432513

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 186750
$ws.Range("J20").Value = 186750
$ws.Range("L20").Value = 186750
$ws.Range("N20").Value = -187222
$ws.Range("H30").Value = 186750
$ws.Range("J30").Value = 186750
$ws.Range("L30").Value = 186750
$ws.Range("N30").Value = -186932
$ws.Range("H62").Value = 2970.8667
$ws.Range("J62").Value = 2953.5715
$ws.Range("L62").Value = 2953.5715
$ws.Range("N62").Value = -4201.5715
$ws.Range("H65").Value = 2970.8667
$ws.Range("J65").Value = 2953.5715
$ws.Range("L65").Value = 14767.8575
$ws.Range("N65").Value = -21007.8575
$ws.Range("H99").Value = 2871
$ws.Range("J99").Value = 2002
$ws.Range("L99").Value = 2002
$ws.Range("N99").Value = -4998
$ws.Range("H107").Value = 1612.7222
$ws.Range("I107").Value = 1029.5
$ws.Range("J107").Value = 1685.625
$ws.Range("K107").Value = 1029.5
$ws.Range("L107").Value = 1685.625
$ws.Range("M107").Value = 890.5
$ws.Range("N107").Value = -5525.625
$ws.Range("H126").Value = 2871
$ws.Range("J126").Value = 2002
$ws.Range("L126").Value = 6006
$ws.Range("N126").Value = -10946
$ws.Range("H128").Value = 186750
$ws.Range("J128").Value = 186750
$ws.Range("L128").Value = 186750
$ws.Range("N128").Value = -196710

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1136.4706
$ws.Range("J34").Value = 1339.8
$ws.Range("L34").Value = 4019.4
$ws.Range("N34").Value = -4187.4
$ws.Range("H39").Value = 4963.375
$ws.Range("J39").Value = 4963.375
$ws.Range("L39").Value = 14890.125
$ws.Range("N39").Value = -15478.125
$ws.Range("H55").Value = 4380
$ws.Range("J55").Value = 5035
$ws.Range("L55").Value = 15105
$ws.Range("N55").Value = -15459
$ws.Range("H109").Value = 3293.3794
$ws.Range("I109").Value = 1773.4546
$ws.Range("J109").Value = 4222.222
$ws.Range("K109").Value = 5320.3638
$ws.Range("L109").Value = 12666.666
$ws.Range("M109").Value = -4280.3638
$ws.Range("N109").Value = -14746.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1528
$ws.Range("I31").Value = 1449.3334
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1449.3334
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -1157.3334
$ws.Range("N31").Value = -2584
$ws.Range("H37").Value = 1528
$ws.Range("I37").Value = 1449.3334
$ws.Range("J37").Value = 2000
$ws.Range("K37").Value = 1449.3334
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = -1172.3334
$ws.Range("N37").Value = -2554
$ws.Range("H70").Value = 24579.867
$ws.Range("I70").Value = 80174.5
$ws.Range("K70").Value = 80174.5
$ws.Range("M70").Value = -79904.5
$ws.Range("H73").Value = 24579.867
$ws.Range("I73").Value = 80174.5
$ws.Range("K73").Value = 80174.5
$ws.Range("M73").Value = -79238.5
$ws.Range("H80").Value = 3576.5
$ws.Range("I80").Value = 3548.1667
$ws.Range("K80").Value = 3548.1667
$ws.Range("M80").Value = -2550.1667
$ws.Range("H83").Value = 3576.5
$ws.Range("I83").Value = 3548.1667
$ws.Range("K83").Value = 17740.8335
$ws.Range("M83").Value = -12748.8335
$ws.Range("H113").Value = 10522.546
$ws.Range("I113").Value = 1443.6666
$ws.Range("J113").Value = 21417.2
$ws.Range("K113").Value = 1443.6666
$ws.Range("L113").Value = 21417.2
$ws.Range("M113").Value = 726.3334
$ws.Range("N113").Value = -25757.2
$ws.Range("H126").Value = 3983.5
$ws.Range("I126").Value = 3943.5
$ws.Range("J126").Value = 4003.5
$ws.Range("K126").Value = 11830.5
$ws.Range("L126").Value = 12010.5
$ws.Range("M126").Value = -9360.5
$ws.Range("N126").Value = -16950.5
$ws.Range("H132").Value = 3095.1428
$ws.Range("I132").Value = 2819.3333
$ws.Range("K132").Value = 8457.999899999999
$ws.Range("M132").Value = -5927.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8659.305
$ws.Range("J46").Value = 8871.046
$ws.Range("L46").Value = 8871.046
$ws.Range("N46").Value = -9247.046
$ws.Range("H132").Value = 6202.615
$ws.Range("I132").Value = 1145
$ws.Range("J132").Value = 6624.0835
$ws.Range("K132").Value = 3435
$ws.Range("L132").Value = 19872.2505
$ws.Range("M132").Value = -905
$ws.Range("N132").Value = -24932.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 113333
$ws.Range("J93").Value = 113333
$ws.Range("L93").Value = 113333
$ws.Range("N93").Value = -118325
$ws.Range("H107").Value = 383.33334
$ws.Range("I107").Value = 316.33334
$ws.Range("J107").Value = 483.83334
$ws.Range("K107").Value = 949.0000200000001
$ws.Range("L107").Value = 1451.50002
$ws.Range("M107").Value = 970.9999799999999
$ws.Range("N107").Value = -5291.500019999999
$ws.Range("H122").Value = 62506944
$ws.Range("I122").Value = 66673670
$ws.Range("K122").Value = 200021010
$ws.Range("M122").Value = -200018560
$ws.Range("H125").Value = 59998.332
$ws.Range("J125").Value = 59998.332
$ws.Range("L125").Value = 59998.332
$ws.Range("N125").Value = -69838.33199999999
$ws.Range("H132").Value = 2759.7307
$ws.Range("I132").Value = 2888.2
$ws.Range("K132").Value = 8664.599999999999
$ws.Range("M132").Value = -6134.599999999999
